$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Paragraph 1 is "Markdown is just a plain text format ..."
$para1 = $tr.Paragraphs(1)

# Insert a new paragraph right after paragraph 1 (before "For example:").
# Leading CR splits the text into a new paragraph, inheriting paragraph 1's
# formatting (no bullet, lvl 0, marL 0) and a plain (no rPr overrides) run.
$newParaText = [char]13 + "It just requires a little training but then it is very fast to edit."
[void]$para1.InsertAfter($newParaText)
